# Daily attendance processing - 2025-10-13 08:28:29
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-ordered "Recorded By" e-mail lists (same sets, new order) ---
$ws.Range("G3").Value  = "Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"
$ws.Range("G4").Value  = "hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg"
$ws.Range("G12").Value = "wessam.atef@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
$ws.Range("G25").Value = "Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"
$ws.Range("G26").Value = "hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg"
$ws.Range("G34").Value = "wessam.atef@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
$ws.Range("G41").Value = "Monica.Eshak@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, marina_atef@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg"

# --- Overall Class Statistics block (K3:L10) ---
$ws.Range("L6").Value  = 8        # Recorded Sessions
$ws.Range("L8").Value  = 31       # Pending Sessions
$ws.Range("L9").Value  = "18.2%"  # Coverage %
$ws.Range("L10").Value = "37.5%"  # Average Attendance %

# --- Per-group statistics row for Year 3 / C2 (row 16) ---
$ws.Range("O16").Value = 5        # Recorded
$ws.Range("Q16").Value = 14       # Pending
$ws.Range("R16").Value = "22.7%"  # Coverage %
$ws.Range("S16").Value = "30.2%"  # Avg Attendance %

# --- Row 42: PHYSIOLOGY session 2 (13/10/2025) now recorded ---
# Copy the "Recorded" (green) formatting from row 41 onto row 42.
$ws.Range("A41:I41").Copy()
$ws.Range("A42:I42").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("G42").Value = "ola.m.abdelfattah@med.asu.edu.eg, marina_atef@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg"
$ws.Range("H42").Value = "15/246"
$ws.Range("I42").Value = "Recorded"
